$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-12-06"

# Update the December row label to reflect the new "through" date
$ws.Range("A13").Value = "December (through 12-06)"

# Update December figures (row 13) for the columns that changed
$ws.Range("C13").Value = 19
$ws.Range("D13").Value = 23
$ws.Range("F13").Value = 6
$ws.Range("G13").Value = 31
$ws.Range("H13").Value = 49

# Update Total row (row 14) for the columns that changed
$ws.Range("C14").Value = 582
$ws.Range("D14").Value = 844
$ws.Range("F14").Value = 540
$ws.Range("G14").Value = 1295
$ws.Range("H14").Value = 1692
